$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A147").Value = "2023-12-09 14:17:39"
$ws.Range("B147").Value = 0.0006000000000000001

$ws.Range("A148").Value = "2023-12-09 14:17:46"
$ws.Range("B148").Value = 0.0004
